# Applies the diff:
#  - Swap match data (columns F:V) between rows 42 and 43
#  - Swap match data (columns F:V) between rows 99 and 100
#  - Append a new row 206 (U. De Chile vs Everton) after the last data row (205)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Swap rows 42 <-> 43 (columns F:V); A:D/A index stay attached to the row.
# ---------------------------------------------------------------------------
$ws.Range("F42").Value = "Union La Calera"
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = "Nublense"
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 2.4
$ws.Range("K42").Value = "20/02/2023 22:12"
$ws.Range("L42").Value = 2
$ws.Range("M42").Value = "26/02/2023 21:57"
$ws.Range("N42").Value = 3.28
$ws.Range("O42").Value = "20/02/2023 22:12"
$ws.Range("P42").Value = 3.6
$ws.Range("Q42").Value = "26/02/2023 21:57"
$ws.Range("R42").Value = 2.95
$ws.Range("S42").Value = "20/02/2023 22:12"
$ws.Range("T42").Value = 3.89
$ws.Range("U42").Value = "26/02/2023 21:53"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/chile/primera-division/union-la-calera-nublense/rXo8WlVS/"

$ws.Range("F43").Value = "Magallanes"
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = "Cobresal"
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 2.74
$ws.Range("K43").Value = "20/02/2023 16:42"
$ws.Range("L43").Value = 2.67
$ws.Range("M43").Value = "26/02/2023 21:57"
$ws.Range("N43").Value = 3.54
$ws.Range("O43").Value = "20/02/2023 16:42"
$ws.Range("P43").Value = 3.55
$ws.Range("Q43").Value = "26/02/2023 21:57"
$ws.Range("R43").Value = 2.59
$ws.Range("S43").Value = "20/02/2023 16:42"
$ws.Range("T43").Value = 2.66
$ws.Range("U43").Value = "26/02/2023 21:58"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/chile/primera-division/magallanes-cobresal/IP0rxg1i/"

# ---------------------------------------------------------------------------
# Swap rows 99 <-> 100 (columns F:V)
# ---------------------------------------------------------------------------
$ws.Range("F99").Value = "Coquimbo"
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = "U. De Chile"
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 2.6
$ws.Range("K99").Value = "03/05/2023 12:43"
$ws.Range("L99").Value = 2.48
$ws.Range("M99").Value = "10/05/2023 23:54"
$ws.Range("N99").Value = 3.29
$ws.Range("O99").Value = "03/05/2023 12:43"
$ws.Range("P99").Value = 3.22
$ws.Range("Q99").Value = "10/05/2023 23:54"
$ws.Range("R99").Value = 2.7
$ws.Range("S99").Value = "03/05/2023 12:43"
$ws.Range("T99").Value = 3.14
$ws.Range("U99").Value = "10/05/2023 23:54"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/chile/primera-division/coquimbo-u-de-chile/l8Zu4N7b/"

$ws.Range("F100").Value = "Cobresal"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "O'Higgins"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1.94
$ws.Range("K100").Value = "03/05/2023 12:43"
$ws.Range("L100").Value = 1.89
$ws.Range("M100").Value = "10/05/2023 23:51"
$ws.Range("N100").Value = 3.61
$ws.Range("O100").Value = "03/05/2023 12:43"
$ws.Range("P100").Value = 3.67
$ws.Range("Q100").Value = "10/05/2023 23:51"
$ws.Range("R100").Value = 3.97
$ws.Range("S100").Value = "03/05/2023 12:43"
$ws.Range("T100").Value = 4.26
$ws.Range("U100").Value = "10/05/2023 23:55"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/chile/primera-division/cobresal-o-higgins/0CG1aaMT/"

# ---------------------------------------------------------------------------
# Append new row 206 (U. De Chile 1 x 2 Everton), mirroring row 205's format.
# ---------------------------------------------------------------------------

# D holds the season as text ("2023"); force text format BEFORE assigning the
# numeric-looking string so the engine doesn't coerce it to a number.
$ws.Range("D206").NumberFormat = "@"
$ws.Range("D206").Value = "2023"

# Copy formatting only (styles: bold/border/center for col A, datetime numFmt
# for col E, default elsewhere) from the previous last row, without touching
# the text we just wrote into D206.
$ws.Range("A205:V205").Copy()
$ws.Range("A206").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A206").Value = 205
$ws.Range("B206").Value = "chile"
$ws.Range("C206").Value = "primera-division"
$ws.Range("E206").Value = 45236.95833333334
$ws.Range("F206").Value = "U. De Chile"
$ws.Range("G206").Value = 1
$ws.Range("H206").Value = "Everton"
$ws.Range("I206").Value = 2
$ws.Range("J206").Value = 2.2
$ws.Range("K206").Value = "03/10/2023 00:12"
$ws.Range("L206").Value = 2.17
$ws.Range("M206").Value = "06/11/2023 22:55"
$ws.Range("N206").Value = 3.4
$ws.Range("O206").Value = "03/10/2023 00:12"
$ws.Range("P206").Value = 3.41
$ws.Range("Q206").Value = "06/11/2023 22:55"
$ws.Range("R206").Value = 3.21
$ws.Range("S206").Value = "03/10/2023 00:12"
$ws.Range("T206").Value = 3.58
$ws.Range("U206").Value = "06/11/2023 22:55"
$ws.Range("V206").Value = "https://www.betexplorer.com/football/chile/primera-division/u-de-chile-everton/IsdExcWD/"
